$wb = $excel.ActiveWorkbook

$wsAbout = $wb.Worksheets.Item("About")
$wsData  = $wb.Worksheets.Item("Data")
$wsTCC   = $wb.Worksheets.Item("TCCpUCD")

# ---------------------------------------------------------------------------
# Data sheet: replace the old "Average cost" row (A13:B13) with a new
# MW-mile / USD / $-per-MW-mile derivation in A1:B3.
# ---------------------------------------------------------------------------
$wsData.Range("A13:B13").ClearContents()

$wsData.Range("A1").Formula = "=17*10^6"
$wsData.Range("B1").Value = "MW-miles"

$wsData.Range("A2").Formula = "=33*10^9"
$wsData.Range("B2").Value = "USD"

$wsData.Range("A3").Formula = "=A2/A1"
$wsData.Range("B3").Value = "$ / MW-mile"

# Remove the embedded picture that no longer applies to the new source.
if ($wsData.Shapes.Count -gt 0) {
    for ($i = $wsData.Shapes.Count; $i -ge 1; $i--) {
        $wsData.Shapes.Item($i).Delete()
    }
}

$wsData.Range("A6").Select()

# ---------------------------------------------------------------------------
# About sheet: update the source citation to the new report.
# ---------------------------------------------------------------------------
$wsAbout.Range("B3").Value = "Americans for a Clean Energy Grid and Grid Strategies"
$wsAbout.Range("B4").Value = 2021
$wsAbout.Range("B5").Value = "Transmission Projects Ready To Go: Plugging Into America's Untapped Renewable Resources"
$wsAbout.Range("B6").Value = "https://cleanenergygrid.org/wp-content/uploads/2019/04/Transmission-Projects-Ready-to-Go-Final.pdf"
$wsAbout.Range("B7").Value = "Pages 11-12"

$wsAbout.Range("A10").Value = "We adjust 2021 dollars to 2012 dollars using the following conversion factor:"
$wsAbout.Range("A11").Value = 0.84730412960844359

# ---------------------------------------------------------------------------
# TCCpUCD sheet: point the cost-per-unit-capacity-distance formula at the
# new Data!A3 ($/MW-mile) figure instead of the old Data!B13 average.
# ---------------------------------------------------------------------------
$wsTCC.Range("B2").Formula = "=Data!A3*About!A11"

# Restore the About sheet as the active tab/selection, matching the
# original session state (only the Data-sheet selection actually moved).
$wsAbout.Activate()
$wsAbout.Range("A11").Select()
